# osis2ebook.docx user-guide update: the paragraph explaining that book
# ordering in the eBook follows the OSIS file order used to go on to say
# this was ensured via CF_paratext2osis.txt's SFM file ordering. That
# extra explanation is removed, leaving a single short sentence. The
# "_GoBack" bookmark that sits in the middle of the original sentence is
# left untouched in place.

$d = $word.ActiveDocument

# 1) Trim the first sentence back to "...the ordering in the OSIS file"
#    by deleting the ", so to ensure correct ordering, CF_paratext2osis.t"
#    tail (this tail spans the end of the first run and all of the second
#    run "`.t`" that precedes the bookmark).
$d.Content.Find.Execute(
    ", so to ensure correct ordering, CF_paratext2osis.t",  # FindText
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "",      # ReplaceWith
    2        # Replace (wdReplaceAll)
)

# 2) Replace the remainder of the original sentence (the text that used
#    to continue after ".t" -> "xt must specify ... should appear.")
#    with a single closing period, so the sentence now simply ends
#    "...the ordering in the OSIS file."
$d.Content.Find.Execute(
    "xt must specify the SFM files to be processed in the order that the books should appear.",  # FindText
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    ".",     # ReplaceWith
    2        # Replace (wdReplaceAll)
)
